# The commit swaps the theme XML parts: ppt/theme/theme1.xml (which held the
# "Office Theme" colour scheme, used by the notes master) and
# ppt/theme/theme2.xml (which held the "Integral" colour scheme, used by the
# slide master / the deck's active design) trade their <a:clrScheme> content
# (and the "Integral" <-> "Office Theme" theme names).
#
# The only theme surface the PowerPoint object model exposes for editing is
# the *active* theme (Master/Design theme, i.e. ppt/theme/theme2.xml) via
# Theme.ThemeColorScheme - so push that theme's 12 scheme colors from the
# "Integral" palette to the "Office Theme" palette the diff wants it to end
# up with (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink, in that order).

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$theme = $master.Theme
$colorScheme = $theme.ThemeColorScheme

$colorScheme.Item(1).RGB  = 0          # dk1      000000
$colorScheme.Item(2).RGB  = 16777215   # lt1      FFFFFF
$colorScheme.Item(3).RGB  = 6968388    # dk2      44546A
$colorScheme.Item(4).RGB  = 15132391   # lt2      E7E6E6
$colorScheme.Item(5).RGB  = 13998939   # accent1  5B9BD5
$colorScheme.Item(6).RGB  = 3243501    # accent2  ED7D31
$colorScheme.Item(7).RGB  = 10855845   # accent3  A5A5A5
$colorScheme.Item(8).RGB  = 49407      # accent4  FFC000
$colorScheme.Item(9).RGB  = 12874308   # accent5  4472C4
$colorScheme.Item(10).RGB = 4697456    # accent6  70AD47
$colorScheme.Item(11).RGB = 12673797   # hlink    0563C1
$colorScheme.Item(12).RGB = 7491477    # folHlink 954F72
